# Beforeschool for all mun (06-22)
# Fill in the "Образование" (Education) category block that previously had
# empty feature-name cells (C24:C25), and move the active selection to D30
# to match where editing left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C24").Value = "Образование"
$ws.Range("C25").Value = "Число мест в дошкол. - beforeschool (шт.) (8014002)"

$ws.Range("D30").Select()
